$wb = $excel.ActiveWorkbook

# --- optimization_parameters sheet: insert the new "L_curve" row and
#     rename the "Model" label to "production_function" ---
$wsOpt = $wb.Worksheets.Item("optimization_parameters")

# Insert a new row above row 9 (estimate_params), shifting everything
# below down by one. The inherited formatting (style 8 on column B)
# matches what the target row needs.
$wsOpt.Rows.Item(9).Insert()

# New row 9: L_curve parameter, default value 0.
$wsOpt.Range("A9").Value = "L_curve"
$wsOpt.Range("B9").Value = 0

# Row 8 label changes from "Model" to "production_function".
$wsOpt.Range("A8").Value = "production_function"

# --- Make optimization_parameters the active sheet / active cell,
#     matching the saved selection in the edited workbook ---
$wsOpt.Activate()
$wsOpt.Range("A9").Select()
